$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title paragraph.
# ---------------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Right before the closing paragraph (the old "Create a feature image for
#    ..." image-prompt paragraph), insert two new paragraphs:
#      a) a bold paragraph repeating the page title
#      b) an italic paragraph holding the (former) meta-description copy
#    then delete the old image-prompt paragraph entirely.
# ---------------------------------------------------------------------------
$oldLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $oldLastPara.Range.Start

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 100 Zombies Free Slot - Review of Features and Payouts</w:t></w:r></w:p>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the thrill of playing 100 Zombies slot for free. Read our review on the features and payouts of the game. Play for free now.</w:t></w:r></w:p>
<w:p><w:r></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.InsertXML($newParasXml)

# The old image-prompt paragraph is now the very last paragraph of the
# document (pushed down by the insert above) - remove it.
$staleLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$staleLastPara.Range.Delete()
